$wb = $excel.ActiveWorkbook

# The source workbook now recalculates automatically instead of manually
$excel.Calculation = -4105

# Add the new "Queries" worksheet at the end of the workbook
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "Queries"

# --- Write cell VALUES first, in the exact order the strings were first
# introduced in the target workbook, so the shared-strings table indices line up. ---

# Row 1 headers that introduce brand-new shared strings
$newSheet.Range("F1").Value = "Query"

$q1 = @'
SELECT CONVERT(decimal(10,2), (100 * (sum(isnull(CallsHandledWithinSLAThreshold,0)))/
 (CAST(ISNULL((CASE WHEN SUM(isnull(PassedCalls,0))+sum(isnull(CallsAbandonedAfterSLAThreshold,0)) = 0 THEN 1
 else SUM(PassedCalls)+sum(isnull(CallsAbandonedAfterSLAThreshold,0)) end),1) AS float)))) AS [Service Level],
 SkillName as [Skill Name],sum([FlowIn]) AS [Flow In],sum([FlowOut]) AS [Flow Out],[SkillId] as [Skill ID],
 [dbo].[SECONDSTOhhmmss](sum(TotalStaffedTIme)/nullif(sum(TotalStaffedAgents),0)) AS [Avg Staff Time],
 [dbo].[SECONDSTOhhmmss](sum(TotalAbandTime)/nullif((sum(AbandCalls)+sum([AcdCalls])),0)) AS [Avg Aband Time],
 sum([AbandCalls]) AS [Aband Calls],
 [dbo].[SECONDSTOhhmmss](sum(SpeedOfAnswer)/nullif(sum([AcdCalls]),0)) AS [Avg Speed Answer],
 [dbo].[SECONDSTOhhmmss](sum([TotalAfterCallTime])) AS [Total After Call Time],
 [dbo].[SECONDSTOhhmmss](sum(TotalTalkTime)/nullif(sum([AcdCalls]),0)) AS [Avg Talk Time],
 sum([AcdCalls])  AS [Total Interaction],
 [dbo].[SECONDSTOhhmmss](sum(TotalAuxTime)) AS [Total Aux Time]
 from [OCM_SkillHistoricalReport]  WITH (NOLOCK)
 WHERE [ReportDateTime]>='ReportBeforeDate' and [ReportDateTime]<='ReportAfterDate' 
 GROUP BY [SkillId],[SkillName]
 ORDER BY [SkillName]
'@
$newSheet.Range("F2").Value = $q1

$newSheet.Range("G1").Value = "QueryDrillGridOne"

$q2 = @'
SELECT CONVERT(decimal(10,2), (100 * (sum(isnull(CallsHandledWithinSLAThreshold,0)))/
(CAST(ISNULL((CASE WHEN SUM(isnull(PassedCalls,0))+sum(isnull(CallsAbandonedAfterSLAThreshold,0)) = 0 THEN 1
else SUM(PassedCalls)+sum(isnull(CallsAbandonedAfterSLAThreshold,0)) end),1) AS float))))  AS [Service Level],
sum([FlowIn]) AS [Flow In],sum([FlowOut]) AS [Flow Out],
[dbo].[SECONDSTOhhmmss](sum(TotalAbandTime)/nullif((sum(AbandCalls)+sum([AcdCalls])),0)) AS [Avg Aband Time],
sum([AbandCalls]) AS [Aband Calls],
Dateint AS [Date],[dbo].[SECONDSTOhhmmss](sum(SpeedOfAnswer)/nullif(sum([AcdCalls]),0)) AS [Avg Speed Answer], 
[dbo].[SECONDSTOhhmmss](sum([TotalAfterCallTime])) AS [Total After Call Time],
[dbo].[SECONDSTOhhmmss](sum(TotalTalkTime)/nullif(sum([AcdCalls]),0)) AS [Avg Talk Time],
sum([AcdCalls]) AS [Total Interaction],
[dbo].[SECONDSTOhhmmss](sum(TotalStaffedTIme)/nullif(sum(TotalStaffedAgents),0)) AS [Avg Staff],
[dbo].[SECONDSTOhhmmss](sum(TotalAuxTime)) AS [Total Aux Time]
from [OCM_SkillHistoricalReport] WITH (NOLOCK)
WHERE [ReportDateTime]>='ReportBeforeDate' and [ReportDateTime]<='ReportAfterDate' and  SkillId like 'SkillIdCapturedFromUI'
GROUP BY [Dateint],[SkillId],[SkillName] ORDER BY [Dateint] ASC
'@
$newSheet.Range("G2").Value = $q2

$newSheet.Range("H1").Value = "QueryDrillGridTwo"

$q3 = @'
SELECT [ServiceLevel] AS [Service Level],[FlowIn] AS [Flow In],[FlowOut] AS [Flow Out],
[dbo].[SECONDSTOhhmmss](TotalStaffedTIme/nullif(TotalStaffedAgents,0)) AS [Avg Staff Time],
[dbo].[SECONDSTOhhmmss](TotalAbandTime/nullif((AbandCalls+[AcdCalls]),0)) AS [Avg Aband Time],
[AbandCalls] AS [Aband Calls],
[dbo].[SECONDSTOhhmmss](SpeedOfAnswer/nullif([AcdCalls],0)) AS [Avg Speed Answer],
[dbo].[SECONDSTOhhmmss]([TotalAfterCallTime]) AS [Total After Call Time],
[dbo].[SECONDSTOhhmmss](TotalTalkTime/nullif([AcdCalls],0)) AS [Avg Talk Time],
[AcdCalls] AS [Total Interaction],[dbo].[SECONDSTOhhmmss](TotalAuxTime) AS [Total Aux Time],[Interval]
FROM [OCM_SkillHistoricalReport] WITH (NOLOCK) WHERE [ReportDateTime]>='ReportBeforeDate' AND [ReportDateTime]<='ReportAfterDate' AND 
[SkillId] like 'SkillIdCapturedFromUI' and [Interval] like '%' 
ORDER BY [intvl] ASC
'@
$newSheet.Range("H2").Value = $q3

# Dates entered as text (leading apostrophe -> quote-prefixed literal text)
$newSheet.Range("D2").Value = "'08-04-2020 00:00:00"
$newSheet.Range("E2").Value = "'22-04-2020 00:00:00"

# Remaining row 1 / row 2 cells -- all reuse already-existing shared strings
$newSheet.Range("A1").Value = "Report Channel"
$newSheet.Range("B1").Value = "Report Name"
$newSheet.Range("C1").Value = "Report Type"
$newSheet.Range("D1").Value = "Start Date"
$newSheet.Range("E1").Value = "End Date"
$newSheet.Range("A2").Value = "Agent"
$newSheet.Range("B2").Value = "OCM Skill Historical Report"
$newSheet.Range("C2").Value = "Date Range"

# --- Styling -------------------------------------------------------------
# Text format ("@") on the header row cells that carry it (all but F1)
$newSheet.Range("A1:E1").NumberFormat = "@"
$newSheet.Range("G1:H1").NumberFormat = "@"

# Wrap the long query text in F2:H2
$newSheet.Range("F2:H2").WrapText = $true

# Row 2 grows to the Excel maximum row height because of the wrapped text
$newSheet.Rows.Item(2).RowHeight = 409.5

# --- Column widths (approximate best-fit widths used by the source file) -
$newSheet.Columns.Item(1).ColumnWidth = 12.98
$newSheet.Columns.Item(2).ColumnWidth = 24.17
$newSheet.Columns.Item(3).ColumnWidth = 10.07
$newSheet.Range("D1:E1").ColumnWidth = 16.98
$newSheet.Columns.Item(6).ColumnWidth = 19.62
$newSheet.Columns.Item(7).ColumnWidth = 15.53
$newSheet.Columns.Item(8).ColumnWidth = 15.26

# --- Selection / active view ----------------------------------------------
$newSheet.Range("E2").Select()
$newSheet.Activate()

Write-Output "Queries sheet added"
